$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Split the "Context" paragraph: trim the trailing sentence about vehicle
#    insurance use-case off of it, and turn that sentence (reworded) plus a
#    new "Real World Use" heading into two new paragraphs right after it.
# ---------------------------------------------------------------------------
$old1 = "for the purpose of data analysis. An example of a visualization of this data would be vehicle specific which could potentially be used by car insurance agencies. "
$new1 = "for the purpose of data analysis. ^pReal World Use ^pAn example of a visualization of this data would be vehicle specific which could potentially be used by car insurance agencies in help determine rates for a location. "
$r1 = $d.Content
$r1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# Bold just the "Real World Use" heading text (not the trailing space).
$r2 = $d.Content
$r2.Find.Execute("Real World Use", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Bold = 1

# ---------------------------------------------------------------------------
# 2. Remove the empty paragraph that sat between the "Data Source" hyperlink
#    paragraph and the "Tools/Technology" heading.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Tools/Technology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$toolsPara = $r3.Paragraphs(1)
$prevPara = $toolsPara.Previous()
$prevPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3. Drop the _GoBack bookmark that used to sit at the end of the
#    "Tools/Technology" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4. Replace the empty paragraph before "GitHub" with a new "Structure"
#    heading paragraph and a paragraph describing the data pipeline.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("GitHub", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$githubPara = $r4.Paragraphs(1)
$structurePara = $githubPara.Previous()
$structurePara.Range.Text = "Structure `rA python program will take the raw data and transform it into usable data for the purpose of this project. This process will filter out unwanted data where the new data will be stored locally. The python program will consist of visualization tools which will produce diagrams and maps that a user may access through a hub like interface for the purpose of analyses. "

# Bold the whole "Structure " paragraph (including its paragraph mark).
$r5 = $d.Content
$r5.Find.Execute("Structure ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r5.Paragraphs(1).Range.Bold = 1

# ---------------------------------------------------------------------------
# 5. Remove the empty paragraph that sat between the "This repo consists..."
#    paragraph and the "Deliverables Order" heading.
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("Deliverables", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deliverablesPara = $r6.Paragraphs(1)
$prevPara2 = $deliverablesPara.Previous()
$prevPara2.Range.Delete()

# ---------------------------------------------------------------------------
# 6. Re-add the _GoBack bookmark, now collapsed at the very start of the
#    "Deliverables Order" paragraph.
# ---------------------------------------------------------------------------
$r7 = $d.Content
$r7.Find.Execute("Deliverables", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deliverablesPara2 = $r7.Paragraphs(1)
$startPos = $deliverablesPara2.Range.Start
$bmRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 7. Strip the trailing run of empty paragraphs at the end of the document.
#    The very last paragraph mark of a document's main story can't itself be
#    deleted (same restriction real Word has), so instead we repeatedly
#    delete the *mark* that ends the "criteria" paragraph -- that merges the
#    next (empty) paragraph into it while leaving its text untouched, which
#    "absorbs" one trailing empty paragraph per iteration.
# ---------------------------------------------------------------------------
$r8 = $d.Content
$r8.Find.Execute("data specific criteria", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$criteriaPara = $r8.Paragraphs(1)
$guard = 0
while ($guard -lt 50) {
    $nextP = $criteriaPara.Next()
    if ($nextP -eq $null -or $nextP.Range.Text -ne "`r") { break }
    $markPos = $criteriaPara.Range.End - 1
    $markRange = $d.Range($markPos, $markPos + 1)
    $markRange.Delete()
    $guard = $guard + 1
}

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
